$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new sample data row (row 2) for dataElement / optionSet / options
$ws.Range("A2").Value = "pGeBz8X2jRq"
$ws.Range("B2").Value = "yHSAPCLxecr"
$ws.Range("C2").Value = "WTSe3FmRFmD"
$ws.Range("D2").Value = 202007
$ws.Range("E2").Value = 8

# Carry over the same number formatting used by the isoPeriod/dataValue header cells
$ws.Range("D1:E1").Copy()
$ws.Range("D2:E2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-apply the values (PasteSpecial of formats only shouldn't disturb them, but make sure)
$ws.Range("D2").Value = 202007
$ws.Range("E2").Value = 8

$ws.Range("A2").Select()
